$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1) - the report grew from 12 columns (A:L) to 37 columns (A:AK)
$headers = @(
  "Title","URL","Status",
  "Title_1","Status_1",
  "Title_2","URL_1","Status_2",
  "Title_1_1","Status_1_1",
  "Title_2_1","URL_1_1","Status_2_1",
  "Title_1_1_1","Status_1_1_1",
  "Title_2_1_1","URL_1_1_1","Status_2_1_1",
  "Title_1_1_1_1","Status_1_1_1_1",
  "Title_2_1_1_1","URL_1_1_1_1","Status_2_1_1_1",
  "Title_1_1_1_1_1","Status_1_1_1_1_1",
  "Title_2_1_1_1_1","URL_1_1_1_1_1","Status_2_1_1_1_1",
  "Title_1_1_1_1_1_1","Status_1_1_1_1_1_1",
  "Title_2_1_1_1_1_1","URL_1_1_1_1_1_1","Status_2_1_1_1_1_1",
  "Title_1_1_1_1_1_1_1","Status_1_1_1_1_1_1_1",
  "Title_2_1_1_1_1_1_1","Status_2_1_1_1_1_1_1"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Latest run's result, appended to the front columns (A:E), replacing
# the previous run's cached values
$ws.Cells.Item(2, 1).Value = "Trackmania"
$ws.Cells.Item(2, 2).Value = "https://www.epicgames.com/store/en-US/product/trackmania/home"
$ws.Cells.Item(2, 3).Value = "Success"
$ws.Cells.Item(2, 4).Value = "Trackmania"
$ws.Cells.Item(2, 5).Value = "Received"
